$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.179.70'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -4.26%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.387.98'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -6.08%  '

$ws.Range('E4').Value = '  +0.39%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '180.72'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -11.89%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '526.07'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.85%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.607'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.47%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.381.75'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -6.12%  '

$ws.Range('E9').Value = '  +0.03%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.622'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.67%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '57.08'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -6.15%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.133'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -12.06%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000253'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -11.54%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.27'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -7.67%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.956.49'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.55%  '

$ws.Range('E16').Value = '  -2.74%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.414.69'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -5.28%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '65.151.78'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.04%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.43'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -7.88%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.14'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -9.95%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.971'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -9.73%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '372.96'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -7.32%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '82.30'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.61%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.73'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -10.43%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.75'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -18.39%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.50'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -8.77%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.64'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -10.00%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.56'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -11.06%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.45'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -10.18%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '675.68'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.20%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '29.56'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.59%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.66'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -20.01%  '

$ws.Range('B33').Value = 'OKB'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '61.52'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.78%  '

$ws.Range('B34').Value = 'Cosmos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.10'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -9.25%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.105'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.70%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.07%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '36.42'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -14.18%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.383'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -9.31%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.00'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.37%  '

$ws.Range('E40').Value = '  -6.68%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.865.12'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -13.14%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.75'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -13.57%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.67'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.80%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0₃0620'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -19.24%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0388'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.36%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.32'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -16.55%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.125'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.96%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '134.91'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.57%  '

$ws.Range('E49').Value = '  -9.08%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.53'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.35%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.59'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -14.39%  '
